$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-15 Thursday", 2) | Out-Null
$d.Content.Find.Execute("88+9=97", $true, $false, $false, $false, $false, $true, 1, $false, "87-31=56", 2) | Out-Null
$d.Content.Find.Execute("57-24=33", $true, $false, $false, $false, $false, $true, 1, $false, "90+4=94", 2) | Out-Null
$d.Content.Find.Execute("97-83=14", $true, $false, $false, $false, $false, $true, 1, $false, "17+63=80", 2) | Out-Null
$d.Content.Find.Execute("10+19=29", $true, $false, $false, $false, $false, $true, 1, $false, "66-55=11", 2) | Out-Null
$d.Content.Find.Execute("11+30=41", $true, $false, $false, $false, $false, $true, 1, $false, "0+59=59", 2) | Out-Null
$d.Content.Find.Execute("22-21=1", $true, $false, $false, $false, $false, $true, 1, $false, "77-75=2", 2) | Out-Null
$d.Content.Find.Execute("59-1=58", $true, $false, $false, $false, $false, $true, 1, $false, "93-80=13", 2) | Out-Null
$d.Content.Find.Execute("83-49=34", $true, $false, $false, $false, $false, $true, 1, $false, "19+60=79", 2) | Out-Null
$d.Content.Find.Execute("1+0=1", $true, $false, $false, $false, $false, $true, 1, $false, "47+44=91", 2) | Out-Null
$d.Content.Find.Execute("2+87=89", $true, $false, $false, $false, $false, $true, 1, $false, "85-54=31", 2) | Out-Null
$d.Content.Find.Execute("61+0=61", $true, $false, $false, $false, $false, $true, 1, $false, "48+37=85", 2) | Out-Null
$d.Content.Find.Execute("88-75=13", $true, $false, $false, $false, $false, $true, 1, $false, "57-44=13", 2) | Out-Null
$d.Content.Find.Execute("77-27=50", $true, $false, $false, $false, $false, $true, 1, $false, "68-49=19", 2) | Out-Null
$d.Content.Find.Execute("79+13=92", $true, $false, $false, $false, $false, $true, 1, $false, "22+17=39", 2) | Out-Null
$d.Content.Find.Execute("14+56=70", $true, $false, $false, $false, $false, $true, 1, $false, "93-81=12", 2) | Out-Null
$d.Content.Find.Execute("58+39=97", $true, $false, $false, $false, $false, $true, 1, $false, "88-43=45", 2) | Out-Null
$d.Content.Find.Execute("23-21=2", $true, $false, $false, $false, $false, $true, 1, $false, "82-67=15", 2) | Out-Null
$d.Content.Find.Execute("13+24=37", $true, $false, $false, $false, $false, $true, 1, $false, "23+9=32", 2) | Out-Null
$d.Content.Find.Execute("97-72=25", $true, $false, $false, $false, $false, $true, 1, $false, "2+62=64", 2) | Out-Null
$d.Content.Find.Execute("67-3=64", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=11", 2) | Out-Null
$d.Content.Find.Execute("68+12=80", $true, $false, $false, $false, $false, $true, 1, $false, "63-40=23", 2) | Out-Null
$d.Content.Find.Execute("50+27=77", $true, $false, $false, $false, $false, $true, 1, $false, "11+37=48", 2) | Out-Null
$d.Content.Find.Execute("32+59=91", $true, $false, $false, $false, $false, $true, 1, $false, "93-40=53", 2) | Out-Null
$d.Content.Find.Execute("15-0=15", $true, $false, $false, $false, $false, $true, 1, $false, "97-16=81", 2) | Out-Null
$d.Content.Find.Execute("49+19=68", $true, $false, $false, $false, $false, $true, 1, $false, "4+15=19", 2) | Out-Null
$d.Content.Find.Execute("25+65=90", $true, $false, $false, $false, $false, $true, 1, $false, "32+46=78", 2) | Out-Null
$d.Content.Find.Execute("68-67=1", $true, $false, $false, $false, $false, $true, 1, $false, "84+8=92", 2) | Out-Null
$d.Content.Find.Execute("70+25=95", $true, $false, $false, $false, $false, $true, 1, $false, "16+43=59", 2) | Out-Null
$d.Content.Find.Execute("79-11=68", $true, $false, $false, $false, $false, $true, 1, $false, "84-36=48", 2) | Out-Null
$d.Content.Find.Execute("74-73=1", $true, $false, $false, $false, $false, $true, 1, $false, "0+67=67", 2) | Out-Null
$d.Content.Find.Execute("27-9=18", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2) | Out-Null
$d.Content.Find.Execute("44+28=72", $true, $false, $false, $false, $false, $true, 1, $false, "73-57=16", 2) | Out-Null
$d.Content.Find.Execute("86-26=60", $true, $false, $false, $false, $false, $true, 1, $false, "47+51=98", 2) | Out-Null
$d.Content.Find.Execute("26+69=95", $true, $false, $false, $false, $false, $true, 1, $false, "27+41=68", 2) | Out-Null
$d.Content.Find.Execute("3+56=59", $true, $false, $false, $false, $false, $true, 1, $false, "20+28=48", 2) | Out-Null
$d.Content.Find.Execute("49-35=14", $true, $false, $false, $false, $false, $true, 1, $false, "60+12=72", 2) | Out-Null
$d.Content.Find.Execute("47-31=16", $true, $false, $false, $false, $false, $true, 1, $false, "74-25=49", 2) | Out-Null
$d.Content.Find.Execute("78+11=89", $true, $false, $false, $false, $false, $true, 1, $false, "90-4=86", 2) | Out-Null
$d.Content.Find.Execute("80+14=94", $true, $false, $false, $false, $false, $true, 1, $false, "23+36=59", 2) | Out-Null
$d.Content.Find.Execute("36+20=56", $true, $false, $false, $false, $false, $true, 1, $false, "32+46=78", 2) | Out-Null
$d.Content.Find.Execute("76+11=87", $true, $false, $false, $false, $false, $true, 1, $false, "25-8=17", 2) | Out-Null
$d.Content.Find.Execute("73-20=53", $true, $false, $false, $false, $false, $true, 1, $false, "72+13=85", 2) | Out-Null
$d.Content.Find.Execute("90-43=47", $true, $false, $false, $false, $false, $true, 1, $false, "68+28=96", 2) | Out-Null
$d.Content.Find.Execute("99-25=74", $true, $false, $false, $false, $false, $true, 1, $false, "61+30=91", 2) | Out-Null
$d.Content.Find.Execute("33-28=5", $true, $false, $false, $false, $false, $true, 1, $false, "50-32=18", 2) | Out-Null
$d.Content.Find.Execute("64-60=4", $true, $false, $false, $false, $false, $true, 1, $false, "26-26=0", 2) | Out-Null
$d.Content.Find.Execute("87-28=59", $true, $false, $false, $false, $false, $true, 1, $false, "17-9=8", 2) | Out-Null
$d.Content.Find.Execute("20+49=69", $true, $false, $false, $false, $false, $true, 1, $false, "36-26=10", 2) | Out-Null
$d.Content.Find.Execute("11+55=66", $true, $false, $false, $false, $false, $true, 1, $false, "56+12=68", 2) | Out-Null
$d.Content.Find.Execute("80+9=89", $true, $false, $false, $false, $false, $true, 1, $false, "8+64=72", 2) | Out-Null
$d.Content.Find.Execute("24-15=9", $true, $false, $false, $false, $false, $true, 1, $false, "27+5=32", 2) | Out-Null
$d.Content.Find.Execute("0+87=87", $true, $false, $false, $false, $false, $true, 1, $false, "16+49=65", 2) | Out-Null
$d.Content.Find.Execute("62-55=7", $true, $false, $false, $false, $false, $true, 1, $false, "82-59=23", 2) | Out-Null
$d.Content.Find.Execute("15+71=86", $true, $false, $false, $false, $false, $true, 1, $false, "13+41=54", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $false, $false, $false, $false, $true, 1, $false, "90-25=65", 2) | Out-Null
$d.Content.Find.Execute("28-11=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+42=48", 2) | Out-Null
$d.Content.Find.Execute("28+19=47", $true, $false, $false, $false, $false, $true, 1, $false, "86-7=79", 2) | Out-Null
$d.Content.Find.Execute("97-14=83", $true, $false, $false, $false, $false, $true, 1, $false, "89-1=88", 2) | Out-Null
$d.Content.Find.Execute("92-72=20", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("87-17=70", $true, $false, $false, $false, $false, $true, 1, $false, "30+37=67", 2) | Out-Null
$d.Content.Find.Execute("6+58=64", $true, $false, $false, $false, $false, $true, 1, $false, "63-42=21", 2) | Out-Null
$d.Content.Find.Execute("65+29=94", $true, $false, $false, $false, $false, $true, 1, $false, "80-20=60", 2) | Out-Null
$d.Content.Find.Execute("77-38=39", $true, $false, $false, $false, $false, $true, 1, $false, "22+32=54", 2) | Out-Null
$d.Content.Find.Execute("87-59=28", $true, $false, $false, $false, $false, $true, 1, $false, "36-7=29", 2) | Out-Null
$d.Content.Find.Execute("95-45=50", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("29+24=53", $true, $false, $false, $false, $false, $true, 1, $false, "11+65=76", 2) | Out-Null
$d.Content.Find.Execute("63+4=67", $true, $false, $false, $false, $false, $true, 1, $false, "59+18=77", 2) | Out-Null
$d.Content.Find.Execute("69-59=10", $true, $false, $false, $false, $false, $true, 1, $false, "14+74=88", 2) | Out-Null
$d.Content.Find.Execute("72-43=29", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=17", 2) | Out-Null
$d.Content.Find.Execute("29-26=3", $true, $false, $false, $false, $false, $true, 1, $false, "24+72=96", 2) | Out-Null
$d.Content.Find.Execute("79+16=95", $true, $false, $false, $false, $false, $true, 1, $false, "89-45=44", 2) | Out-Null
$d.Content.Find.Execute("57+21=78", $true, $false, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("43+20=63", $true, $false, $false, $false, $false, $true, 1, $false, "3+38=41", 2) | Out-Null
$d.Content.Find.Execute("96-94=2", $true, $false, $false, $false, $false, $true, 1, $false, "39+7=46", 2) | Out-Null
$d.Content.Find.Execute("74-0=74", $true, $false, $false, $false, $false, $true, 1, $false, "3+35=38", 2) | Out-Null
$d.Content.Find.Execute("32+32=64", $true, $false, $false, $false, $false, $true, 1, $false, "7+7=14", 2) | Out-Null
$d.Content.Find.Execute("69+14=83", $true, $false, $false, $false, $false, $true, 1, $false, "69-56=13", 2) | Out-Null
$d.Content.Find.Execute("84-46=38", $true, $false, $false, $false, $false, $true, 1, $false, "52+10=62", 2) | Out-Null
$d.Content.Find.Execute("50-40=10", $true, $false, $false, $false, $false, $true, 1, $false, "46-33=13", 2) | Out-Null
$d.Content.Find.Execute("53-28=25", $true, $false, $false, $false, $false, $true, 1, $false, "60-11=49", 2) | Out-Null
$d.Content.Find.Execute("47+3=50", $true, $false, $false, $false, $false, $true, 1, $false, "87-63=24", 2) | Out-Null
$d.Content.Find.Execute("77-20=57", $true, $false, $false, $false, $false, $true, 1, $false, "8+91=99", 2) | Out-Null
$d.Content.Find.Execute("54+6=60", $true, $false, $false, $false, $false, $true, 1, $false, "58+14=72", 2) | Out-Null
$d.Content.Find.Execute("54-36=18", $true, $false, $false, $false, $false, $true, 1, $false, "38+39=77", 2) | Out-Null
$d.Content.Find.Execute("44+6=50", $true, $false, $false, $false, $false, $true, 1, $false, "33-25=8", 2) | Out-Null
$d.Content.Find.Execute("47-47=0", $true, $false, $false, $false, $false, $true, 1, $false, "37+19=56", 2) | Out-Null
$d.Content.Find.Execute("52+37=89", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=73", 2) | Out-Null
$d.Content.Find.Execute("35+52=87", $true, $false, $false, $false, $false, $true, 1, $false, "47+7=54", 2) | Out-Null
$d.Content.Find.Execute("14+37=51", $true, $false, $false, $false, $false, $true, 1, $false, "32-24=8", 2) | Out-Null
$d.Content.Find.Execute("76-3=73", $true, $false, $false, $false, $false, $true, 1, $false, "17+31=48", 2) | Out-Null
$d.Content.Find.Execute("37+42=79", $true, $false, $false, $false, $false, $true, 1, $false, "95-24=71", 2) | Out-Null
$d.Content.Find.Execute("84-75=9", $true, $false, $false, $false, $false, $true, 1, $false, "16+12=28", 2) | Out-Null
$d.Content.Find.Execute("18+79=97", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=2", 2) | Out-Null
$d.Content.Find.Execute("88-40=48", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=90", 2) | Out-Null
$d.Content.Find.Execute("21+26=47", $true, $false, $false, $false, $false, $true, 1, $false, "37+3=40", 2) | Out-Null
$d.Content.Find.Execute("8+67=75", $true, $false, $false, $false, $false, $true, 1, $false, "61-20=41", 2) | Out-Null
$d.Content.Find.Execute("62-11=51", $true, $false, $false, $false, $false, $true, 1, $false, "25+72=97", 2) | Out-Null
$d.Content.Find.Execute("37+25=62", $true, $false, $false, $false, $false, $true, 1, $false, "50-26=24", 2) | Out-Null
$d.Content.Find.Execute("12+0=12", $true, $false, $false, $false, $false, $true, 1, $false, "4+41=45", 2) | Out-Null
$d.Content.Find.Execute("85-61=24", $true, $false, $false, $false, $false, $true, 1, $false, "29-2=27", 2) | Out-Null
